$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right after "总计" (total) sheet,
#    pushing 2022-Q2 / 2022-Q1 / ... down by one position.
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$ns = $wb.Worksheets.Add($null, $wsTotal)
$ns.Name = "2022-Q3"

# Reference to the (now shifted) former "2022-Q2" sheet, used purely as a
# formatting template so the new sheet's header row / A-column match the
# existing quarterly sheets exactly (bold header style, bordered index col).
$srcTemplate = $wb.Worksheets.Item(3)

# ---- header row (row 1) ----
$srcTemplate.Range("B1:H1").Copy()
$ns.Range("B1:H1").PasteSpecial(-4122)
$ns.Range("B1").Value = "基金代码"
$ns.Range("C1").Value = "基金名称"
$ns.Range("D1").Value = "基金规模"
$ns.Range("E1").Value = "股票总仓位"
$ns.Range("F1").Value = "仓位占比"
$ns.Range("G1").Value = "持有市值(亿元)"
$ns.Range("H1").Value = "仓位排名"

# ---- column A (row index, bold/bordered style like the other sheets) ----
$srcTemplate.Range("A2").Copy()
$ns.Range("A2:A18").PasteSpecial(-4122)

# ---- data rows 2-18 ----
$ns.Range("A2").Value = 0
$ns.Range("B2:G2").NumberFormat = "@"
$ns.Range("B2").Value = "360006"
$ns.Range("C2").Value = "光大保德信新增长混合"
$ns.Range("D2").Value = "17.68"
$ns.Range("E2").Value = "83.30"
$ns.Range("F2").Value = "3.83"
$ns.Range("G2").Value = "0.6771"
$ns.Range("H2").Value = 8
$ns.Range("A3").Value = 1
$ns.Range("B3:G3").NumberFormat = "@"
$ns.Range("B3").Value = "161914"
$ns.Range("C3").Value = "万家创业板2年定期开放混合A"
$ns.Range("D3").Value = "7.30"
$ns.Range("E3").Value = "99.15"
$ns.Range("F3").Value = "8.03"
$ns.Range("G3").Value = "0.5862"
$ns.Range("H3").Value = 5
$ns.Range("A4").Value = 2
$ns.Range("B4:G4").NumberFormat = "@"
$ns.Range("B4").Value = "009837"
$ns.Range("C4").Value = "华夏磐锐一年定期开放混合A"
$ns.Range("D4").Value = "14.02"
$ns.Range("E4").Value = "94.15"
$ns.Range("F4").Value = "3.98"
$ns.Range("G4").Value = "0.5580"
$ns.Range("H4").Value = 9
$ns.Range("A5").Value = 3
$ns.Range("B5:G5").NumberFormat = "@"
$ns.Range("B5").Value = "003713"
$ns.Range("C5").Value = "英大睿盛灵活配置混合A"
$ns.Range("D5").Value = "2.83"
$ns.Range("E5").Value = "93.65"
$ns.Range("F5").Value = "7.16"
$ns.Range("G5").Value = "0.2026"
$ns.Range("H5").Value = 4
$ns.Range("A6").Value = 4
$ns.Range("B6:G6").NumberFormat = "@"
$ns.Range("B6").Value = "003714"
$ns.Range("C6").Value = "英大睿盛灵活配置混合C"
$ns.Range("D6").Value = "2.19"
$ns.Range("E6").Value = "93.65"
$ns.Range("F6").Value = "7.16"
$ns.Range("G6").Value = "0.1568"
$ns.Range("H6").Value = 4
$ns.Range("A7").Value = 5
$ns.Range("B7:G7").NumberFormat = "@"
$ns.Range("B7").Value = "010676"
$ns.Range("C7").Value = "光大保德信新机遇混合"
$ns.Range("D7").Value = "2.85"
$ns.Range("E7").Value = "84.08"
$ns.Range("F7").Value = "5.40"
$ns.Range("G7").Value = "0.1539"
$ns.Range("H7").Value = 5
$ns.Range("A8").Value = 6
$ns.Range("B8:G8").NumberFormat = "@"
$ns.Range("B8").Value = "161915"
$ns.Range("C8").Value = "万家创业板2年定期开放混合C"
$ns.Range("D8").Value = "1.43"
$ns.Range("E8").Value = "99.15"
$ns.Range("F8").Value = "8.03"
$ns.Range("G8").Value = "0.1148"
$ns.Range("H8").Value = 5
$ns.Range("A9").Value = 7
$ns.Range("B9:G9").NumberFormat = "@"
$ns.Range("B9").Value = "001607"
$ns.Range("C9").Value = "英大策略优选混合A"
$ns.Range("D9").Value = "0.57"
$ns.Range("E9").Value = "91.98"
$ns.Range("F9").Value = "6.35"
$ns.Range("G9").Value = "0.0362"
$ns.Range("H9").Value = 5
$ns.Range("A10").Value = 8
$ns.Range("B10:G10").NumberFormat = "@"
$ns.Range("B10").Value = "012522"
$ns.Range("C10").Value = "英大稳固增强核心一年持有混合C"
$ns.Range("D10").Value = "1.24"
$ns.Range("E10").Value = "27.71"
$ns.Range("F10").Value = "1.59"
$ns.Range("G10").Value = "0.0197"
$ns.Range("H10").Value = 5
$ns.Range("A11").Value = 9
$ns.Range("B11:G11").NumberFormat = "@"
$ns.Range("B11").Value = "003447"
$ns.Range("C11").Value = "英大睿鑫灵活配置混合C"
$ns.Range("D11").Value = "0.21"
$ns.Range("E11").Value = "92.71"
$ns.Range("F11").Value = "7.99"
$ns.Range("G11").Value = "0.0168"
$ns.Range("H11").Value = 2
$ns.Range("A12").Value = 10
$ns.Range("B12:G12").NumberFormat = "@"
$ns.Range("B12").Value = "009838"
$ns.Range("C12").Value = "华夏磐锐一年定期开放混合C"
$ns.Range("D12").Value = "0.39"
$ns.Range("E12").Value = "94.15"
$ns.Range("F12").Value = "3.98"
$ns.Range("G12").Value = "0.0155"
$ns.Range("H12").Value = 9
$ns.Range("A13").Value = 11
$ns.Range("B13:G13").NumberFormat = "@"
$ns.Range("B13").Value = "007152"
$ns.Range("C13").Value = "诺德策略精选混合"
$ns.Range("D13").Value = "0.34"
$ns.Range("E13").Value = "80.74"
$ns.Range("F13").Value = "4.08"
$ns.Range("G13").Value = "0.0139"
$ns.Range("H13").Value = 8
$ns.Range("A14").Value = 12
$ns.Range("B14:G14").NumberFormat = "@"
$ns.Range("B14").Value = "012521"
$ns.Range("C14").Value = "英大稳固增强核心一年持有混合A"
$ns.Range("D14").Value = "0.75"
$ns.Range("E14").Value = "27.71"
$ns.Range("F14").Value = "1.59"
$ns.Range("G14").Value = "0.0119"
$ns.Range("H14").Value = 5
$ns.Range("A15").Value = 13
$ns.Range("B15:G15").NumberFormat = "@"
$ns.Range("B15").Value = "007133"
$ns.Range("C15").Value = "嘉实长青竞争优势股票A"
$ns.Range("D15").Value = "0.24"
$ns.Range("E15").Value = "90.21"
$ns.Range("F15").Value = "4.88"
$ns.Range("G15").Value = "0.0117"
$ns.Range("H15").Value = 10
$ns.Range("A16").Value = 14
$ns.Range("B16:G16").NumberFormat = "@"
$ns.Range("B16").Value = "003446"
$ns.Range("C16").Value = "英大睿鑫灵活配置混合A"
$ns.Range("D16").Value = "0.07"
$ns.Range("E16").Value = "92.71"
$ns.Range("F16").Value = "7.99"
$ns.Range("G16").Value = "0.0056"
$ns.Range("H16").Value = 2
$ns.Range("A17").Value = 15
$ns.Range("B17:G17").NumberFormat = "@"
$ns.Range("B17").Value = "007134"
$ns.Range("C17").Value = "嘉实长青竞争优势股票C"
$ns.Range("D17").Value = "0.04"
$ns.Range("E17").Value = "90.21"
$ns.Range("F17").Value = "4.88"
$ns.Range("G17").Value = "0.0020"
$ns.Range("H17").Value = 10
$ns.Range("A18").Value = 16
$ns.Range("B18:G18").NumberFormat = "@"
$ns.Range("B18").Value = "001608"
$ns.Range("C18").Value = "英大策略优选混合C"
$ns.Range("D18").Value = "0.02"
$ns.Range("E18").Value = "91.98"
$ns.Range("F18").Value = "6.35"
$ns.Range("G18").Value = "0.0013"
$ns.Range("H18").Value = 5

# ------------------------------------------------------------------
# 2. Update the "总计" (total) sheet: insert a new row 2 for 2022-Q3,
#    shifting the existing quarters down by one row.
# ------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 17
$wsTotal.Range("D2").Value = 2.58

# The row-index column (A) is a simple 0-based counter; since a new row0
# was inserted, every pre-existing row's counter shifts up by one too.
for ($r = 3; $r -le 8; $r++) {
    $cur = $wsTotal.Range("A$r").Value2
    $wsTotal.Range("A$r").Value = $cur + 1
}
